# AllInOneGenerator: add Protocol.proto automation support.
# Reorders the KDAInfo/PlayerInfo/SC_ON_ACCEPT block (rows 2-9) so the
# nested-message types (KDAInfo, PlayerInfo) are defined before the
# message that references them (SC_ON_ACCEPT), and appends two new
# messages (SC_SEND_MESSAGE / CS_SEND_MESSAGE) at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder rows 2-9 -------------------------------------------------
# Row 2: KDAInfo / kill / uint32
$ws.Range("A2").Value = "KDAInfo"
$ws.Range("B2").Value = "kill"
$ws.Range("C2").Value = "uint32"
$ws.Range("D2").Value = $null

# Row 3: KDAInfo / death / uint32
$ws.Range("A3").Value = "KDAInfo"
$ws.Range("B3").Value = "death"
$ws.Range("C3").Value = "uint32"

# Row 4: KDAInfo / assist / uint32
$ws.Range("A4").Value = "KDAInfo"
$ws.Range("B4").Value = "assist"
$ws.Range("C4").Value = "uint32"

# Row 5: now blank (spacer row)
$ws.Range("A5").Value = $null
$ws.Range("B5").Value = $null
$ws.Range("C5").Value = $null
$ws.Range("D5").Value = $null

# Row 6: PlayerInfo / playerId / uint32
$ws.Range("A6").Value = "PlayerInfo"
$ws.Range("B6").Value = "playerId"
$ws.Range("C6").Value = "uint32"

# Row 7: PlayerInfo / kda / KDAInfo
$ws.Range("A7").Value = "PlayerInfo"
$ws.Range("B7").Value = "kda"
$ws.Range("C7").Value = "KDAInfo"

# Row 8: now blank (spacer row)
$ws.Range("A8").Value = $null
$ws.Range("B8").Value = $null
$ws.Range("C8").Value = $null
$ws.Range("D8").Value = $null

# Row 9: SC_ON_ACCEPT / playerId / uint32 / 서버로 부터 부여받은 ID
$ws.Range("A9").Value = "SC_ON_ACCEPT"
$ws.Range("B9").Value = "playerId"
$ws.Range("C9").Value = "uint32"
$ws.Range("D9").Value = "서버로 부터 부여받은 ID"

# --- Append new SC_SEND_MESSAGE / CS_SEND_MESSAGE messages -------------
# Row 121-122: SC_SEND_MESSAGE
$ws.Range("A121").Value = "SC_SEND_MESSAGE"
$ws.Range("B121").Value = "playerId"
$ws.Range("C121").Value = "uint32"
$ws.Range("D121").Value = "플레이어의 고유 ID"

$ws.Range("A122").Value = "SC_SEND_MESSAGE"
$ws.Range("B122").Value = "message"
$ws.Range("C122").Value = "string"
$ws.Range("D122").Value = "전송할 문자열"

# Row 124-125: CS_SEND_MESSAGE
$ws.Range("A124").Value = "CS_SEND_MESSAGE"
$ws.Range("B124").Value = "playerId"
$ws.Range("C124").Value = "uint32"
$ws.Range("D124").Value = "플레이어의 고유 ID"

$ws.Range("A125").Value = "CS_SEND_MESSAGE"
$ws.Range("B125").Value = "message"
$ws.Range("C125").Value = "string"
$ws.Range("D125").Value = "전송할 문자열"

# --- Update the view so the newly added rows are visible --------------
$ws.Range("B129").Select()
$excel.ActiveWindow.ScrollRow = 106
